$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) values: force-keep as text via NumberFormat toggle,
# since many look like plain numbers and would otherwise be auto-coerced
# to numeric values by the COM "smart" Value assignment.
$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value2 = '64.480.26'
$cell.ClearFormats()
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value2 = '3.155.95'
$cell.ClearFormats()
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value2 = '607.37'
$cell.ClearFormats()
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value2 = '146.74'
$cell.ClearFormats()
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value2 = '3.149.63'
$cell.ClearFormats()
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value2 = '0.475'
$cell.ClearFormats()
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value2 = '35.99'
$cell.ClearFormats()
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value2 = '3.673.21'
$cell.ClearFormats()
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value2 = '64.422.94'
$cell.ClearFormats()
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value2 = '3.153.98'
$cell.ClearFormats()
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value2 = '6.94'
$cell.ClearFormats()
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value2 = '478.74'
$cell.ClearFormats()
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value2 = '14.57'
$cell.ClearFormats()
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value2 = '7.71'
$cell.ClearFormats()
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value2 = '13.76'
$cell.ClearFormats()
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value2 = '83.40'
$cell.ClearFormats()
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value2 = '2.88'
$cell.ClearFormats()
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value2 = '6.79'
$cell.ClearFormats()
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value2 = '0.114'
$cell.ClearFormats()
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value2 = '26.17'
$cell.ClearFormats()
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value2 = '1.11'
$cell.ClearFormats()
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value2 = '6.02'
$cell.ClearFormats()
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value2 = '54.15'
$cell.ClearFormats()
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value2 = '0.0₃0717'
$cell.ClearFormats()
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value2 = '451.89'
$cell.ClearFormats()
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value2 = '0.119'
$cell.ClearFormats()
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value2 = '8.44'
$cell.ClearFormats()
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value2 = '2.842.89'
$cell.ClearFormats()
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value2 = '26.37'
$cell.ClearFormats()
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value2 = '118.74'
$cell.ClearFormats()

# Other text columns (Coin name / Link / Volume%) assign directly.
$ws.Range('E2').Value2 = '  -3.20%  '
$ws.Range('E3').Value2 = '  -2.62%  '
$ws.Range('E4').Value2 = '  +0.11%  '
$ws.Range('E5').Value2 = '  +0.49%  '
$ws.Range('E6').Value2 = '  -6.59%  '
$ws.Range('E7').Value2 = '  +0.06%  '
$ws.Range('E8').Value2 = '  -2.83%  '
$ws.Range('E9').Value2 = '  -3.61%  '
$ws.Range('E10').Value2 = '  -7.35%  '
$ws.Range('E11').Value2 = '  -4.16%  '
$ws.Range('E12').Value2 = '  -5.50%  '
$ws.Range('E13').Value2 = '  -6.15%  '
$ws.Range('E14').Value2 = '  -8.01%  '
$ws.Range('E15').Value2 = '  -2.65%  '
$ws.Range('E16').Value2 = '  -3.35%  '
$ws.Range('E17').Value2 = '  +1.04%  '
$ws.Range('E18').Value2 = '  -2.67%  '
$ws.Range('E19').Value2 = '  -4.92%  '
$ws.Range('E20').Value2 = '  -5.93%  '
$ws.Range('E21').Value2 = '  -4.77%  '
$ws.Range('E22').Value2 = '  -5.06%  '
$ws.Range('E23').Value2 = '  -4.05%  '
$ws.Range('E24').Value2 = '  -6.12%  '
$ws.Range('E25').Value2 = '  -3.11%  '
$ws.Range('E27').Value2 = '  -4.02%  '
$ws.Range('E28').Value2 = '  -7.75%  '
$ws.Range('E29').Value2 = '  -6.98%  '
$ws.Range('B30').Value2 = 'NEARProtocol'
$ws.Range('C30').Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E30').Value2 = '  -1.90%  '
$ws.Range('B31').Value2 = 'Hedera'
$ws.Range('C31').Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E31').Value2 = '  -32.13%  '
$ws.Range('E32').Value2 = '  -5.79%  '
$ws.Range('E33').Value2 = '  +0.06%  '
$ws.Range('E35').Value2 = '  -4.62%  '
$ws.Range('E36').Value2 = '  -5.50%  '
$ws.Range('E37').Value2 = '  -2.05%  '
$ws.Range('E38').Value2 = '  -10.97%  '
$ws.Range('E39').Value2 = '  -8.72%  '
$ws.Range('E40').Value2 = '  -10.42%  '
$ws.Range('E41').Value2 = '  -6.18%  '
$ws.Range('B42').Value2 = 'Kaspa'
$ws.Range('C42').Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E42').Value2 = '  -7.36%  '
$ws.Range('B43').Value2 = 'Cosmos'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E43').Value2 = '  -3.47%  '
$ws.Range('E44').Value2 = '  -3.66%  '
$ws.Range('E45').Value2 = '  -8.32%  '
$ws.Range('E46').Value2 = '  -8.40%  '
$ws.Range('E47').Value2 = '  -6.67%  '
$ws.Range('E48').Value2 = '  -0.05%  '
$ws.Range('E49').Value2 = '  -4.12%  '
$ws.Range('E50').Value2 = '  -4.35%  '
$ws.Range('E51').Value2 = '  -1.88%  '
